$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 45229
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 2000
$ws.Range("P3").Value = 667

# Row 4
$ws.Range("D4").Value = 45203
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 2500
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2500
$ws.Range("P4").Value = 833

# Row 5
$ws.Range("D5").Value = 45203
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1500
$ws.Range("P5").Value = 500

# Row 6
$ws.Range("D6").Value = 45205
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 400
$ws.Range("L6").Value = 2500
$ws.Range("M6").Value = 2250
$ws.Range("P6").Value = 750

# Row 7
$ws.Range("D7").Value = 44838
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 1200
$ws.Range("L7").Value = 1300
$ws.Range("M7").Value = 1250
$ws.Range("P7").Value = 417

# Row 8
$ws.Range("D8").Value = 44838
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 1000
$ws.Range("P8").Value = 333

# Row 9
$ws.Range("D9").Value = 45148
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 2500
$ws.Range("L9").Value = 2500
$ws.Range("M9").Value = 2500
$ws.Range("P9").Value = 833

# Row 10
$ws.Range("D10").Value = 45148
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 2000
$ws.Range("P10").Value = 667

# Row 11
$ws.Range("D11").Value = 45145
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = 2500
$ws.Range("M11").Value = 2500
$ws.Range("P11").Value = 833

# Row 12
$ws.Range("D12").Value = 45145
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = 2000
$ws.Range("P12").Value = 667

# Row 13
$ws.Range("D13").Value = 45135
$ws.Range("J13").Value = 70

# Row 14
$ws.Range("D14").Value = 45161
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2500
$ws.Range("P14").Value = 833

# Row 15
$ws.Range("D15").Value = 45191
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 2500
$ws.Range("L15").Value = 2500
$ws.Range("M15").Value = 2500
$ws.Range("P15").Value = 833

# Row 16
$ws.Range("D16").Value = 45175
$ws.Range("J16").Value = 150

# Row 17
$ws.Range("D17").Value = 45246
$ws.Range("J17").Value = 300

# Row 18
$ws.Range("D18").Value = 45246
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = 2000
$ws.Range("P18").Value = 667

# Row 19
$ws.Range("D19").Value = 45215
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 2000
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = 2000
$ws.Range("P19").Value = 667

# Row 20
$ws.Range("D20").Value = 45163
$ws.Range("J20").Value = 100

# Row 21
$ws.Range("D21").Value = 45247
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 2500
$ws.Range("L21").Value = 2500
$ws.Range("M21").Value = 2500
$ws.Range("P21").Value = 833

# Row 22
$ws.Range("D22").Value = 45247
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 2000
$ws.Range("P22").Value = 667

# Row 23
$ws.Range("D23").Value = 45134
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 2500
$ws.Range("L23").Value = 2500
$ws.Range("M23").Value = 2500
$ws.Range("P23").Value = 833

# Row 24
$ws.Range("D24").Value = 45133
$ws.Range("J24").Value = 80

# Row 27
$ws.Range("D27").Value = 45225
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = 2000
$ws.Range("P27").Value = 667

# Row 28
$ws.Range("D28").Value = 45149
$ws.Range("J28").Value = 80
$ws.Range("K28").Value = 2500
$ws.Range("L28").Value = 2500
$ws.Range("M28").Value = 2500
$ws.Range("P28").Value = 833

# Row 29
$ws.Range("D29").Value = 45149
$ws.Range("I29").Value = "Segunda"
$ws.Range("J29").Value = 80
$ws.Range("K29").Value = 2000
$ws.Range("L29").Value = 2000
$ws.Range("M29").Value = 2000
$ws.Range("P29").Value = 667

# Row 30
$ws.Range("D30").Value = 45233

# Row 31
$ws.Range("D31").Value = 45233
$ws.Range("I31").Value = "Segunda"
$ws.Range("J31").Value = 250
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = 1500
$ws.Range("P31").Value = 500

# Row 32
$ws.Range("D32").Value = 44846
$ws.Range("J32").Value = 200
$ws.Range("K32").Value = 1200
$ws.Range("L32").Value = 1300
$ws.Range("M32").Value = 1250
$ws.Range("P32").Value = 417

# Row 33
$ws.Range("D33").Value = 44846
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 1000
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = 1000
$ws.Range("P33").Value = 333

# Row 34
$ws.Range("D34").Value = 45217
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 400
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = 2250
$ws.Range("P34").Value = 750

# Row 35
$ws.Range("D35").Value = 45176

# Row 36
$ws.Range("D36").Value = 44832
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 1200
$ws.Range("L36").Value = 1300
$ws.Range("M36").Value = 1250
$ws.Range("P36").Value = 417

# Row 37
$ws.Range("D37").Value = 44832
$ws.Range("I37").Value = "Segunda"
$ws.Range("J37").Value = 150
$ws.Range("K37").Value = 1000
$ws.Range("L37").Value = 1000
$ws.Range("M37").Value = 1000
$ws.Range("P37").Value = 333

# Row 38
$ws.Range("D38").Value = 45244
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 2500
$ws.Range("L38").Value = 2500
$ws.Range("M38").Value = 2500
$ws.Range("P38").Value = 833

# Row 39
$ws.Range("D39").Value = 45244
$ws.Range("I39").Value = "Segunda"
$ws.Range("J39").Value = 100
$ws.Range("K39").Value = 2000
$ws.Range("L39").Value = 2000
$ws.Range("M39").Value = 2000
$ws.Range("P39").Value = 667

# Row 40
$ws.Range("D40").Value = 45195
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = 2500
$ws.Range("P40").Value = 833

# Row 41
$ws.Range("D41").Value = 45219
$ws.Range("J41").Value = 300
$ws.Range("K41").Value = 2000
$ws.Range("L41").Value = 2500
$ws.Range("M41").Value = 2250
$ws.Range("P41").Value = 750

# Row 42
$ws.Range("D42").Value = 45166
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 120
$ws.Range("K42").Value = 2500
$ws.Range("L42").Value = 2500
$ws.Range("M42").Value = 2500
$ws.Range("P42").Value = 833

# Row 43
$ws.Range("D43").Value = 45173
$ws.Range("J43").Value = 100
